# Financials update: insert a new "latest period" column before column D on the
# FNRN worksheet, shifting the existing D:K data right to E:L, then fill in the
# new column D with the newest period's figures (a new fiscal year-end date in
# row 7/38/80, and the corresponding financial figures in every data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D. This shifts the old D:K -> E:L automatically,
#    including values, and extends dimension/row spans.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D has no number formatting / style yet.
#    Copy the formats from column E (which now holds what used to be column D)
#    into column D so dates keep the date format and values keep the
#    thousands-separator number format used throughout the sheet.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the newest period's values, row by row.
#    Empty-string entries ("NA") reuse the sheet's existing "NA" shared text.
$newValues = @{
    7  = 43465
    8  = 45600
    9  = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 3400
    18 = 42200
    20 = -25000
    21 = 17900
    22 = 0
    23 = 17300
    24 = 4700
    25 = 0
    26 = 12600
    27 = 12600
    28 = 0
    29 = "NA"
    30 = 0
    31 = 0
    32 = 25000
    33 = 12600
    34 = 0
    35 = 12600
    38 = 43465
    41 = 116000
    42 = 13600
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 6600
    49 = 1600
    50 = 0
    51 = 0
    52 = 6000
    53 = 0
    54 = 1249800
    57 = 12800
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 1137400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 23900
    73 = 0
    74 = 0
    75 = 0
    76 = 112500
    77 = 0
    80 = 43465
    81 = 12600
    83 = 600
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 15700
    91 = -1000
    92 = 0
    93 = 0
    94 = -72500
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 20000
    101 = 0
    102 = -36900
}

foreach ($r in $newValues.Keys) {
    $addr = "D" + $r
    $ws.Range($addr).Value2 = $newValues[$r]
}

Write-Output "done"
